$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire row 7 (CUSTOM CLUB 6 / STEVE), shifting rows 8-20 up by one.
$ws.Rows.Item(7).Delete()

# Restore the selection to match the recorded state after the edit.
$ws.Range("C19").Select()
